# CDS Data Validations script updates
# Rename the "WebExcel" tab/column references to "ExDataExcel" / "...ExcelData.xlsx"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "ExDataExcel"
$ws.Range("D2").Value = "TC02_CDSValidation_by_ParticipantID - 2_ExcelData.xlsx"
$ws.Range("D3").Value = "TC02_CDSValidation_by_ParticipantID - 2_ExcelData.xlsx"
$ws.Range("D4").Value = "TC02_CDSValidation_by_ParticipantID - 2_ExcelData.xlsx"

$ws.Range("D2").Select()
